# Adapt column header formatting to respective input file names (#7)
#
# - Rename "<Spaltenname>_old" headers (A1:J1) to "<Spaltenname>_FV2410"
# - Rename "<Spaltenname>_new" headers (L1:U1) to "<Spaltenname>_FV2504"
# - Freeze the header row
# - Turn the data range into a native Excel Table ("Table1") with autofilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename headers -----------------------------------------------------
# Columns A-J use the "_old" suffix, columns L-U use the "_new" suffix; the
# "diff" header in column K is left untouched.
$ws.Range("A1:J1").Replace("_old", "_FV2410", 2, 1, $false, $false, $false) | Out-Null
$ws.Range("L1:U1").Replace("_new", "_FV2504", 2, 1, $false, $false, $false) | Out-Null

# --- 2. Freeze the top (header) row ----------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into an Excel table --------------------------
$dataRange = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$lo.Name = "Table1"
